# Insert a new data row at row 61, shifting existing rows 61..146 down to 62..147.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(61).Insert()

$ws.Range("A61").Value = 8
$ws.Range("B61").Value = "Terminal La Palmera de La Serena"
$ws.Range("C61").Value = "Coquimbo"
$ws.Range("D61").Value = 44540
$ws.Range("E61").Value = 4
$ws.Range("F61").Value = 100112021
$ws.Range("G61").Value = "Ají"
$ws.Range("H61").Value = "Inferno"
$ws.Range("I61").Value = "Primera"
$ws.Range("J61").Value = 300
$ws.Range("K61").Value = 18000
$ws.Range("L61").Value = 19000
$ws.Range("M61").Value = 18500
$ws.Range("N61").Value = "$/caja 15 kilos"
$ws.Range("O61").Value = "Provincia de Limarí"
$ws.Range("P61").Value = 1233
$ws.Range("Q61").Value = 15
$ws.Range("R61").Value = "Hortaliza"
